# Auto-generated script applying value updates per commit diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 102.454544
$ws.Range("I33").Value = 102.454544
$ws.Range("K33").Value = 102.454544
$ws.Range("M33").Value = 126.545456
$ws.Range("H40").Value = 1640
$ws.Range("J40").Value = 2500
$ws.Range("L40").Value = 2500
$ws.Range("N40").Value = -2850
$ws.Range("H74").Value = 5685977
$ws.Range("I74").Value = 3539.8
$ws.Range("J74").Value = 7357282.5
$ws.Range("K74").Value = 3539.8
$ws.Range("L74").Value = 7357282.5
$ws.Range("M74").Value = -2603.8
$ws.Range("N74").Value = -7359154.5
$ws.Range("H77").Value = 5685977
$ws.Range("I77").Value = 3539.8
$ws.Range("J77").Value = 7357282.5
$ws.Range("K77").Value = 17699
$ws.Range("L77").Value = 36786412.5
$ws.Range("M77").Value = -13019
$ws.Range("N77").Value = -36795772.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1569.25
$ws.Range("J2").Value = 1672.5454
$ws.Range("L2").Value = 1672.5454
$ws.Range("N2").Value = -1898.5454
$ws.Range("H32").Value = 21107
$ws.Range("J32").Value = 17657
$ws.Range("L32").Value = 17657
$ws.Range("N32").Value = -18231
$ws.Range("H61").Value = 1607.9362
$ws.Range("I61").Value = 1032.6757
$ws.Range("K61").Value = 1032.6757
$ws.Range("M61").Value = -820.6757
$ws.Range("H116").Value = 1569.25
$ws.Range("J116").Value = 1672.5454
$ws.Range("L116").Value = 1672.5454
$ws.Range("N116").Value = -6260.5454
$ws.Range("H122").Value = 2344.5454
$ws.Range("I122").Value = 2337.1428
$ws.Range("K122").Value = 7011.428400000001
$ws.Range("M122").Value = -4561.428400000001
$ws.Range("H132").Value = 12392.108
$ws.Range("I132").Value = 1403.3784
$ws.Range("J132").Value = 57568
$ws.Range("K132").Value = 4210.135200000001
$ws.Range("L132").Value = 172704
$ws.Range("M132").Value = -1680.135200000001
$ws.Range("N132").Value = -177764
$ws.Range("H136").Value = 1607.9362
$ws.Range("I136").Value = 1032.6757
$ws.Range("K136").Value = 3098.0271
$ws.Range("M136").Value = -548.0271000000002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1569.25
$ws.Range("J3").Value = 1672.5454
$ws.Range("L3").Value = 1672.5454
$ws.Range("N3").Value = -1900.5454
$ws.Range("H22").Value = 328
$ws.Range("I22").Value = 328
$ws.Range("K22").Value = 328
$ws.Range("M22").Value = -155
$ws.Range("H134").Value = 30383.527
$ws.Range("I134").Value = 35073.97
$ws.Range("J134").Value = 1302.8
$ws.Range("K134").Value = 105221.91
$ws.Range("L134").Value = 3908.4
$ws.Range("M134").Value = -102686.91
$ws.Range("N134").Value = -8978.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1083.625
$ws.Range("I16").Value = 1113.3334
$ws.Range("J16").Value = 994.5
$ws.Range("K16").Value = 1113.3334
$ws.Range("L16").Value = 994.5
$ws.Range("M16").Value = -826.3334
$ws.Range("N16").Value = -1568.5
$ws.Range("H31").Value = 18678.389
$ws.Range("I31").Value = 21347.467
$ws.Range("K31").Value = 21347.467
$ws.Range("M31").Value = -21052.467
$ws.Range("H34").Value = 18678.389
$ws.Range("I34").Value = 21347.467
$ws.Range("K34").Value = 21347.467
$ws.Range("M34").Value = -21145.467
$ws.Range("H94").Value = 4372.6665
$ws.Range("I94").Value = 2685.3333
$ws.Range("J94").Value = 6060
$ws.Range("K94").Value = 2685.3333
$ws.Range("L94").Value = 6060
$ws.Range("M94").Value = -2234.3333
$ws.Range("N94").Value = -6962
$ws.Range("H113").Value = 1083.625
$ws.Range("I113").Value = 1113.3334
$ws.Range("J113").Value = 994.5
$ws.Range("K113").Value = 1113.3334
$ws.Range("L113").Value = 994.5
$ws.Range("M113").Value = 1056.6666
$ws.Range("N113").Value = -5334.5
$ws.Range("H132").Value = 33408.234
$ws.Range("I132").Value = 44477.25
$ws.Range("J132").Value = 6842.6
$ws.Range("K132").Value = 133431.75
$ws.Range("L132").Value = 20527.8
$ws.Range("M132").Value = -130901.75
$ws.Range("N132").Value = -25587.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 939.75
$ws.Range("I17").Value = 180
$ws.Range("J17").Value = 1193
$ws.Range("K17").Value = 540
$ws.Range("L17").Value = 3579
$ws.Range("M17").Value = -371
$ws.Range("N17").Value = -3917
$ws.Range("H22").Value = 9058.333000000001
$ws.Range("I22").Value = 11733.333
$ws.Range("K22").Value = 35199.999
$ws.Range("M22").Value = -35030.999
$ws.Range("H27").Value = 9058.333000000001
$ws.Range("I27").Value = 11733.333
$ws.Range("K27").Value = 35199.999
$ws.Range("M27").Value = -35097.999
$ws.Range("H75").Value = 1042
$ws.Range("I75").Value = 1096.3334
$ws.Range("J75").Value = 933.3333
$ws.Range("K75").Value = 3289.0002
$ws.Range("L75").Value = 2799.9999
$ws.Range("M75").Value = -2291.0002
$ws.Range("N75").Value = -4795.9999
$ws.Range("H78").Value = 1042
$ws.Range("I78").Value = 1096.3334
$ws.Range("J78").Value = 933.3333
$ws.Range("K78").Value = 9867.000599999999
$ws.Range("L78").Value = 8399.9997
$ws.Range("M78").Value = -4875.000599999999
$ws.Range("N78").Value = -18383.9997
$ws.Range("H92").Value = 496.66666
$ws.Range("I92").Value = 496.66666
$ws.Range("K92").Value = 1489.99998
$ws.Range("M92").Value = -241.9999800000001
$ws.Range("H131").Value = 777.29
$ws.Range("I131").Value = 564.1667
$ws.Range("J131").Value = 790.8936
$ws.Range("K131").Value = 1692.5001
$ws.Range("L131").Value = 2372.6808
$ws.Range("M131").Value = 3347.4999
$ws.Range("N131").Value = -12452.6808
$ws.Range("H139").Value = 2374.0625
$ws.Range("I139").Value = 1754.7
$ws.Range("K139").Value = 5264.1
$ws.Range("M139").Value = -124.1000000000004

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4360.6787
$ws.Range("I126").Value = 3244.389
$ws.Range("J126").Value = 6370
$ws.Range("K126").Value = 9733.167000000001
$ws.Range("L126").Value = 19110
$ws.Range("M126").Value = -7263.167000000001
$ws.Range("N126").Value = -24050
$ws.Range("H132").Value = 65374.457
$ws.Range("I132").Value = 52948.7
$ws.Range("K132").Value = 158846.1
$ws.Range("M132").Value = -156316.1

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5181.2354
$ws.Range("I7").Value = 5442.2856
$ws.Range("J7").Value = 3963
$ws.Range("K7").Value = 5442.2856
$ws.Range("L7").Value = 3963
$ws.Range("M7").Value = -5330.2856
$ws.Range("N7").Value = -4187
$ws.Range("H126").Value = 5181.2354
$ws.Range("I126").Value = 5442.2856
$ws.Range("J126").Value = 3963
$ws.Range("K126").Value = 16326.8568
$ws.Range("L126").Value = 11889
$ws.Range("M126").Value = -13856.8568
$ws.Range("N126").Value = -16829
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("H113").Value = 1690289.2
$ws.Range("I113").Value = 1345.3636
$ws.Range("K113").Value = 4036.0908
$ws.Range("M113").Value = -1866.0908
$ws.Range("H126").Value = 1557.6471
$ws.Range("I126").Value = 1072.5
$ws.Range("J126").Value = 1706.9231
$ws.Range("K126").Value = 3217.5
$ws.Range("L126").Value = 5120.7693
$ws.Range("M126").Value = -747.5
$ws.Range("N126").Value = -10060.7693
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492
